$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.840.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.975.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.33%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.977.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.34%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.145"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.13%  "

$ws.Range("E12").Value = "  +2.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.463.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.521.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.959.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "444.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.49%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0887"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.72%  "

$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("E40").Value = "  +2.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.119"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.283"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.699.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "369.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.88%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0342"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("E51").Value = "  -0.53%  "
